$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.627.70"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "2.565.20"
$ws.Range("E3").Value = "  -3.60%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'521.20"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -1.35%  "
$ws.Range("D9").Value = "2.577.65"
$ws.Range("E9").Value = "  -3.42%  "
$ws.Range("E10").Value = "  -3.94%  "
$ws.Range("D11").Value = "'0.101"
$ws.Range("E11").Value = "  -2.24%  "
$ws.Range("E12").Value = "  -2.80%  "
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "3.017.21"
$ws.Range("E14").Value = "  -3.58%  "
$ws.Range("D15").Value = "57.609.83"
$ws.Range("E15").Value = "  -2.06%  "
$ws.Range("D16").Value = "'20.18"
$ws.Range("E16").Value = "  -3.93%  "
$ws.Range("E17").Value = "  -2.36%  "
$ws.Range("D18").Value = "2.574.15"
$ws.Range("E18").Value = "  -3.58%  "
$ws.Range("D19").Value = "'335.94"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").Value = "'4.29"
$ws.Range("D21").Value = "'10.20"
$ws.Range("E21").Value = "  -2.00%  "
$ws.Range("D22").Value = "'6.23"
$ws.Range("E22").Value = "  -2.28%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "'65.24"
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").Value = "'0.403"
$ws.Range("E26").Value = "  -4.94%  "
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").Value = "2.679.12"
$ws.Range("E28").Value = "  -3.72%  "
$ws.Range("E29").Value = "  -2.84%  "
$ws.Range("D30").Value = "0.0₃0749"
$ws.Range("E30").Value = "  -6.99%  "
$ws.Range("D32").Value = "'6.25"
$ws.Range("E32").Value = "  -6.55%  "
$ws.Range("D34").Value = "'18.65"
$ws.Range("E34").Value = "  -1.32%  "
$ws.Range("D35").Value = "'148.69"
$ws.Range("E35").Value = "  -1.40%  "
$ws.Range("E36").Value = "  -2.95%  "
$ws.Range("E37").Value = "  -4.02%  "
$ws.Range("D38").Value = "'0.845"
$ws.Range("E38").Value = "  -9.61%  "
$ws.Range("D39").Value = "'36.14"
$ws.Range("E39").Value = "  -1.74%  "
$ws.Range("D40").Value = "'0.830"
$ws.Range("E40").Value = "  -5.32%  "
$ws.Range("E41").Value = "  -1.45%  "
$ws.Range("D42").Value = "'3.52"
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("D44").Value = "'268.89"
$ws.Range("E44").Value = "  -2.71%  "
$ws.Range("D45").Value = "'0.0957"
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("D47").Value = "'0.589"
$ws.Range("E47").Value = "  -3.70%  "
$ws.Range("D48").Value = "'18.87"
$ws.Range("E48").Value = "  -4.30%  "
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("D50").Value = "1.970.83"
$ws.Range("E50").Value = "  -4.45%  "
$ws.Range("E51").Value = "  -2.40%  "
